$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.174425666666667
$ws.Range("H2").Value = 18.523277
$ws.Range("I2").Value = 0.4796039576068858
$ws.Range("J2").Value = 0.5423873011998577
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.02354
$ws.Range("N2").Value = 0.07062
$ws.Range("O2").Value = 0.01363425747161943
$ws.Range("P2").Value = 0.01984867005197698
$ws.Range("Q2").Value = 0.1453459801933334
$ws.Range("R2").Value = 1.30811382174
$ws.Range("S2").Value = 0.006539043842419931
$ws.Range("T2").Value = 0.01076566658189823

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.174425666666667
$ws.Range("H3").Value = 18.523277
$ws.Range("I3").Value = 0.4796039576068858
$ws.Range("J3").Value = 0.5423873011998577
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.03558433333333334
$ws.Range("N3").Value = 0.106753
$ws.Range("O3").Value = 0.02061027878600664
$ws.Range("P3").Value = 0.03000431993852589
$ws.Range("Q3").Value = 0.2197128210645556
$ws.Range("R3").Value = 1.977415389581
$ws.Range("S3").Value = 0.009884771273150027
$ws.Range("T3").Value = 0.01627396211579414

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.174425666666667
$ws.Range("H4").Value = 18.523277
$ws.Range("I4").Value = 0.4796039576068858
$ws.Range("J4").Value = 0.5423873011998577
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04573
$ws.Range("N4").Value = 0.13719
$ws.Range("O4").Value = 0.02648660128195227
$ws.Range("P4").Value = 0.03855903489706489
$ws.Range("Q4").Value = 0.2823564857366667
$ws.Range("R4").Value = 2.54120837163
$ws.Range("S4").Value = 0.01270307879837993
$ws.Range("T4").Value = 0.02091393087469015

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.174425666666667
$ws.Range("H5").Value = 18.523277
$ws.Range("I5").Value = 0.4796039576068858
$ws.Range("J5").Value = 0.5423873011998577
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.621679
$ws.Range("N5").Value = 3.243358
$ws.Range("O5").Value = 0.9392688624604216
$ws.Range("P5").Value = 0.9115879751124322
$ws.Range("Q5").Value = 10.01293644069433
$ws.Range("R5").Value = 60.077618644166
$ws.Range("S5").Value = 0.4504770636929359
$ws.Range("T5").Value = 0.4944337416274751

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.228940333333334
$ws.Range("H6").Value = 6.686821
$ws.Range("I6").Value = 0.1731349056329954
$ws.Range("J6").Value = 0.1957994147469982
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.02354
$ws.Range("N6").Value = 0.07062
$ws.Range("O6").Value = 0.01363425747161943
$ws.Range("P6").Value = 0.01984867005197698
$ws.Range("Q6").Value = 0.05246925544666668
$ws.Range("R6").Value = 0.47222329902
$ws.Range("S6").Value = 0.002360565880724792
$ws.Range("T6").Value = 0.003886357979683364

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.228940333333334
$ws.Range("H7").Value = 6.686821
$ws.Range("I7").Value = 0.1731349056329954
$ws.Range("J7").Value = 0.1957994147469982
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.03558433333333334
$ws.Range("N7").Value = 0.106753
$ws.Range("O7").Value = 0.02061027878600664
$ws.Range("P7").Value = 0.03000431993852589
$ws.Range("Q7").Value = 0.07931535580144446
$ws.Range("R7").Value = 0.713838202213
$ws.Range("S7").Value = 0.003568358672684986
$ws.Range("T7").Value = 0.005874828283845061

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ntrk3"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.228940333333334
$ws.Range("H8").Value = 6.686821
$ws.Range("I8").Value = 0.1731349056329954
$ws.Range("J8").Value = 0.1957994147469982
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04573
$ws.Range("N8").Value = 0.13719
$ws.Range("O8").Value = 0.02648660128195227
$ws.Range("P8").Value = 0.03855903489706489
$ws.Range("Q8").Value = 0.1019294414433333
$ws.Range("R8").Value = 0.9173649729900001
$ws.Range("S8").Value = 0.004585755213489581
$ws.Range("T8").Value = 0.007549836466054387

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ntrk3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.228940333333334
$ws.Range("H9").Value = 6.686821
$ws.Range("I9").Value = 0.1731349056329954
$ws.Range("J9").Value = 0.1957994147469982
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.621679
$ws.Range("N9").Value = 3.243358
$ws.Range("O9").Value = 0.9392688624604216
$ws.Range("P9").Value = 0.9115879751124322
$ws.Range("Q9").Value = 3.614625730819667
$ws.Range("R9").Value = 21.687754384918
$ws.Range("S9").Value = 0.162620225866096
$ws.Range("T9").Value = 0.1784883920174154

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ntf3"
$ws.Range("C10").Value = "Ntrk3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.470643
$ws.Range("H10").Value = 8.941286
$ws.Range("I10").Value = 0.3472611367601188
$ws.Range("J10").Value = 0.2618132840531441
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.02354
$ws.Range("N10").Value = 0.07062
$ws.Range("O10").Value = 0.01363425747161943
$ws.Range("P10").Value = 0.01984867005197698
$ws.Range("Q10").Value = 0.10523893622
$ws.Range("R10").Value = 0.63143361732
$ws.Range("S10").Value = 0.004734647748474706
$ws.Range("T10").Value = 0.005196645490395383

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Ntf3"
$ws.Range("C11").Value = "Ntrk3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.470643
$ws.Range("H11").Value = 8.941286
$ws.Range("I11").Value = 0.3472611367601188
$ws.Range("J11").Value = 0.2618132840531441
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.03558433333333334
$ws.Range("N11").Value = 0.106753
$ws.Range("O11").Value = 0.02061027878600664
$ws.Range("P11").Value = 0.03000431993852589
$ws.Range("Q11").Value = 0.1590848507263334
$ws.Range("R11").Value = 0.954509104358
$ws.Range("S11").Value = 0.007157148840171627
$ws.Range("T11").Value = 0.007855529538886694

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Ntf3"
$ws.Range("C12").Value = "Ntrk3"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.470643
$ws.Range("H12").Value = 8.941286
$ws.Range("I12").Value = 0.3472611367601188
$ws.Range("J12").Value = 0.2618132840531441
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.04573
$ws.Range("N12").Value = 0.13719
$ws.Range("O12").Value = 0.02648660128195227
$ws.Range("P12").Value = 0.03855903489706489
$ws.Range("Q12").Value = 0.20444250439
$ws.Range("R12").Value = 1.22665502634
$ws.Range("S12").Value = 0.009197767270082766
$ws.Range("T12").Value = 0.01009526755632034

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Ntf3"
$ws.Range("C13").Value = "Ntrk3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.470643
$ws.Range("H13").Value = 8.941286
$ws.Range("I13").Value = 0.3472611367601188
$ws.Range("J13").Value = 0.2618132840531441
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.621679
$ws.Range("N13").Value = 3.243358
$ws.Range("O13").Value = 0.9392688624604216
$ws.Range("P13").Value = 0.9115879751124322
$ws.Range("Q13").Value = 7.249947869597
$ws.Range("R13").Value = 28.999791478388
$ws.Range("S13").Value = 0.3261715729013897
$ws.Range("T13").Value = 0.2386658414675416
